$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: two new Betarraga price records (date 2021-09-10 / serial 44449,
# grades "Primera" and "Segunda") are inserted right after existing row 206.
# That pushes all the following weekly records down by two rows; the two
# rows that fall off the bottom of the original range (227/228) land on the
# newly created rows 229/230.
#
# Columns: D=Fecha, I=Calidad, J=Volumen, K=Precio minimo, L=Precio maximo,
#          M=Precio promedio ponderado, P=Precio $/Kg

$data = @(
  @(206, 44449, "Primera", 5200, 90, 100, 95, 95),
  @(207, 44449, "Segunda", 2500, 70, 80, 75, 75),
  @(208, 44161, "Primera", 17000, 80, 100, 89, 89),
  @(209, 44428, "Primera", 5200, 100, 110, 105, 105),
  @(210, 44428, "Segunda", 2500, 90, 90, 90, 90),
  @(211, 44442, "Primera", 5200, 90, 100, 95, 95),
  @(212, 44442, "Segunda", 2500, 70, 80, 75, 75),
  @(213, 44435, "Primera", 27100, 90, 110, 100, 100),
  @(214, 44435, "Segunda", 13600, 90, 100, 81, 81),
  @(215, 44319, "Primera", 2500, 130, 130, 130, 130),
  @(216, 44175, "Primera", 12000, 90, 100, 95, 95),
  @(217, 44376, "Primera", 5200, 100, 110, 105, 105),
  @(218, 44376, "Segunda", 2500, 80, 80, 80, 80),
  @(219, 44412, "Primera", 5200, 100, 120, 95, 95),
  @(220, 44412, "Segunda", 2500, 70, 70, 70, 70),
  @(221, 44223, "Primera", 9000, 100, 120, 111, 111),
  @(222, 44223, "Segunda", 7000, 80, 80, 80, 80),
  @(223, 44314, "Primera", 5200, 130, 130, 130, 130),
  @(224, 44448, "Primera", 5200, 90, 100, 95, 95),
  @(225, 44448, "Segunda", 2500, 70, 80, 75, 75),
  @(226, 44167, "Primera", 15000, 80, 100, 91, 91),
  @(227, 44399, "Primera", 6100, 90, 100, 95, 95),
  @(228, 44399, "Segunda", 2500, 70, 70, 70, 70),
  @(229, 44400, "Primera", 6100, 90, 100, 95, 95),
  @(230, 44400, "Segunda", 2500, 70, 70, 70, 70)
)

$dateFormat = $ws.Range("D206").NumberFormat

foreach ($rec in $data) {
  $r = $rec[0]

  # Newly created rows need the whole record written out (they don't exist
  # in the original sheet yet); existing rows only need the fields that
  # actually change, but writing every field is simplest and harmless.
  if ($r -gt 228) {
    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100114014
    $ws.Cells.Item($r, 7).Value = "Betarraga"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 14).Value = "`$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
  }

  $ws.Cells.Item($r, 4).Value = $rec[1]
  $ws.Cells.Item($r, 9).Value = $rec[2]
  $ws.Cells.Item($r, 10).Value = $rec[3]
  $ws.Cells.Item($r, 11).Value = $rec[4]
  $ws.Cells.Item($r, 12).Value = $rec[5]
  $ws.Cells.Item($r, 13).Value = $rec[6]
  $ws.Cells.Item($r, 16).Value = $rec[7]
}
